$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the scraped cryptos.com snapshot (price + 1h volume change).
# D/E values are forced to Text via a leading apostrophe so numeric-looking
# strings like '1.006' keep matching the source t="inlineStr" cells, then
# the style is reset to 'Normal' so no stray quote-prefix formatting sticks.

$ws.Range('D2').Value = "'27.572.55"
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = "'  +1.15%  "
$ws.Range('E2').Style = "Normal"
$ws.Range('D3').Value = "'1.763.58"
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = "'  -1.20%  "
$ws.Range('E3').Style = "Normal"
$ws.Range('D4').Value = "'1.006"
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = "'  +0.18%  "
$ws.Range('E4').Style = "Normal"
$ws.Range('D5').Value = "'336.08"
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = "'  +0.40%  "
$ws.Range('E5').Style = "Normal"
$ws.Range('D6').Value = "'1.002"
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = "'  +0.04%  "
$ws.Range('E6').Style = "Normal"
$ws.Range('D7').Value = "'0.3834"
$ws.Range('D7').Style = "Normal"
$ws.Range('D8').Value = "'0.3404"
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = "'  -0.71%  "
$ws.Range('E8').Style = "Normal"
$ws.Range('D9').Value = "'46.73"
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = "'  -4.13%  "
$ws.Range('E9').Style = "Normal"
$ws.Range('E10').Value = "'  -4.93%  "
$ws.Range('E10').Style = "Normal"
$ws.Range('D11').Value = "'0.07401"
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = "'  -1.19%  "
$ws.Range('E11').Style = "Normal"
$ws.Range('D12').Value = "'1.003"
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = "'  +0.13%  "
$ws.Range('E12').Style = "Normal"
$ws.Range('D13').Value = "'22.37"
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = "'  +2.46%  "
$ws.Range('E13').Style = "Normal"
$ws.Range('D14').Value = "'6.337"
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = "'  -2.10%  "
$ws.Range('E14').Style = "Normal"
$ws.Range('D15').Value = "'1.761.76"
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = "'  -1.46%  "
$ws.Range('E15').Style = "Normal"
$ws.Range('D16').Value = "'7.025"
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = "'  -0.96%  "
$ws.Range('E16').Style = "Normal"
$ws.Range('D17').Value = "'0.00001073"
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = "'  -2.22%  "
$ws.Range('E17').Style = "Normal"
$ws.Range('D18').Value = "'0.06646"
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = "'  -0.03%  "
$ws.Range('E18').Style = "Normal"
$ws.Range('D19').Value = "'82.14"
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = "'  -2.00%  "
$ws.Range('E19').Style = "Normal"
$ws.Range('D20').Value = "'1.003"
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = "'  +0.14%  "
$ws.Range('E20').Style = "Normal"
$ws.Range('D21').Value = "'17.34"
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = "'  +0.00%  "
$ws.Range('E21').Style = "Normal"
$ws.Range('D22').Value = "'6.397"
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = "'  -3.82%  "
$ws.Range('E22').Style = "Normal"
$ws.Range('D23').Value = "'27.584.90"
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = "'  +1.18%  "
$ws.Range('E23').Style = "Normal"
$ws.Range('E24').Value = "'  -2.32%  "
$ws.Range('E24').Style = "Normal"
$ws.Range('D25').Value = "'2.373"
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = "'  -1.62%  "
$ws.Range('E25').Style = "Normal"
$ws.Range('D26').Value = "'20.66"
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = "'  -3.06%  "
$ws.Range('E26').Style = "Normal"
$ws.Range('B27').Value = "LidoDAOToken"
$ws.Range('C27').Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range('D27').Value = "'2.438"
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = "'  -4.20%  "
$ws.Range('E27').Style = "Normal"
$ws.Range('B28').Value = "ImmutableX"
$ws.Range('C28').Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range('D28').Value = "'1.418"
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = "'  -5.06%  "
$ws.Range('E28').Style = "Normal"
$ws.Range('D29').Value = "'153.14"
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = "'  +0.08%  "
$ws.Range('E29').Style = "Normal"
$ws.Range('D30').Value = "'134.79"
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = "'  +0.62%  "
$ws.Range('E30').Style = "Normal"
$ws.Range('D31').Value = "'1.961.51"
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = "'  -1.41%  "
$ws.Range('E31').Style = "Normal"
$ws.Range('D32').Value = "'6.095"
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = "'  +0.10%  "
$ws.Range('E32').Style = "Normal"
$ws.Range('D33').Value = "'3.959"
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = "'  -1.71%  "
$ws.Range('E33').Style = "Normal"
$ws.Range('D34').Value = "'0.08775"
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = "'  +0.88%  "
$ws.Range('E34').Style = "Normal"
$ws.Range('D35').Value = "'12.70"
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = "'  -3.90%  "
$ws.Range('E35').Style = "Normal"
$ws.Range('D36').Value = "'0.02405"
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = "'  +2.95%  "
$ws.Range('E36').Style = "Normal"
$ws.Range('D37').Value = "'5.352"
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = "'  -1.75%  "
$ws.Range('E37').Style = "Normal"
$ws.Range('D38').Value = "'0.6767"
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = "'  -2.51%  "
$ws.Range('E38').Style = "Normal"
$ws.Range('D39').Value = "'0.2183"
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = "'  -0.97%  "
$ws.Range('E39').Style = "Normal"
$ws.Range('D40').Value = "'0.06267"
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = "'  -0.87%  "
$ws.Range('E40').Style = "Normal"
$ws.Range('B41').Value = "TrustWalletToken"
$ws.Range('C41').Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range('D41').Value = "'1.254"
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = "'  +1.26%  "
$ws.Range('E41').Style = "Normal"
$ws.Range('B42').Value = "WEMIXTOKEN"
$ws.Range('C42').Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range('D42').Value = "'1.505"
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = "'  -9.46%  "
$ws.Range('E42').Style = "Normal"
$ws.Range('D43').Value = "'8.302"
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = "'  -6.00%  "
$ws.Range('E43').Style = "Normal"
$ws.Range('D44').Value = "'14.26"
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = "'  -1.41%  "
$ws.Range('E44').Style = "Normal"
$ws.Range('D45').Value = "'1.002"
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = "'  +0.07%  "
$ws.Range('E45').Style = "Normal"
$ws.Range('D46').Value = "'0.6237"
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = "'  -4.43%  "
$ws.Range('E46').Style = "Normal"
$ws.Range('D47').Value = "'3.827"
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = "'  -0.39%  "
$ws.Range('E47').Style = "Normal"
$ws.Range('D48').Value = "'132.35"
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = "'  +2.47%  "
$ws.Range('E48').Style = "Normal"
$ws.Range('D49').Value = "'2.078"
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = "'  -3.38%  "
$ws.Range('E49').Style = "Normal"
$ws.Range('D50').Value = "'0.07380"
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = "'  +3.46%  "
$ws.Range('E50').Style = "Normal"
$ws.Range('B51').Value = "EOS"
$ws.Range('C51').Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range('D51').Value = "'1.144"
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = "'  +2.51%  "
$ws.Range('E51').Style = "Normal"
